# Fix Training Data Issue (#48)
# The "Date" column (BF) was populated with the literal file-name-derived
# string "6-28-2011-12" on every data row. That text looks like
# "June 28, 2011-12" but the NBA season column actually needs an
# unambiguous ISO 8601 date: 2012-06-28 (the season that started in 2011
# wrapped up in June 2012, and NBA stats for that night were reported a
# day off from the naive "6-28-2011-12" parse).
#
# Cells BF2:BF31 hold that text as plain strings (inlineStr), not real
# dates, and must stay plain text after the edit -- so we force text entry
# with a leading apostrophe (quote prefix) to stop Excel's automatic date
# recognition from turning "2012-06-28" into a date serial number, then
# immediately clear the incidental formatting that the quote prefix
# trick applies, so the cell falls back to its original (unstyled) look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 58)   # column BF = 58
    if ($cell.Text -eq "6-28-2011-12") {
        $cell.Value = "'2012-06-28"
        $cell.ClearFormats()
    }
}
